$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells that are being updated,
# so numeric-looking strings (e.g. "1.000", "237.50") keep their exact
# original text representation instead of being coerced to numbers.
$dCells = @("D2","D3","D4","D5","D6","D8","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.311.60'
$ws.Range("E2").Value = '  -1.23%  '
$ws.Range("D3").Value = '1.879.48'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '237.50'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  -2.36%  '
$ws.Range("D8").Value = '0.2890'
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("E9").Value = '  -2.77%  '
$ws.Range("D10").Value = '1.879.66'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("D11").Value = '16.92'
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("D12").Value = '0.07397'
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = '5.183'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '87.97'
$ws.Range("E14").Value = '  -1.07%  '
$ws.Range("D15").Value = '0.6594'
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("D16").Value = '30.283.06'
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("D17").Value = '13.58'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '0.000007727'
$ws.Range("E19").Value = '  -2.86%  '
$ws.Range("D20").Value = '5.470'
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("D21").Value = '2.136.51'
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").Value = '195.94'
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").Value = '6.152'
$ws.Range("E24").Value = '  -2.52%  '
$ws.Range("D25").Value = '9.423'
$ws.Range("E25").Value = '  -2.44%  '
$ws.Range("D26").Value = '163.34'
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("D27").Value = '18.22'
$ws.Range("E27").Value = '  -3.74%  '
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("D29").Value = '1.440'
$ws.Range("E29").Value = '  -2.83%  '
$ws.Range("D30").Value = '4.271'
$ws.Range("E30").Value = '  -2.27%  '
$ws.Range("D31").Value = '0.09144'
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").Value = '4.042'
$ws.Range("E32").Value = '  -0.56%  '
$ws.Range("D33").Value = '0.05050'
$ws.Range("E33").Value = '  -4.54%  '
$ws.Range("D34").Value = '0.7408'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").Value = '1.138'
$ws.Range("E35").Value = '  +1.71%  '
$ws.Range("D36").Value = '2.708'
$ws.Range("E36").Value = '  -0.89%  '
$ws.Range("D37").Value = '0.01836'
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("E38").Value = '  -3.39%  '
$ws.Range("D39").Value = '0.9144'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("E40").Value = '  -0.59%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '5.876'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4319'
$ws.Range("E43").Value = '  -3.16%  '
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").Value = '7.625'
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = '0.1347'
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("D47").Value = '1.565'
$ws.Range("E47").Value = '  +8.72%  '
$ws.Range("D48").Value = '64.97'
$ws.Range("E48").Value = '  -12.62%  '
$ws.Range("D49").Value = '8.845'
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("D50").Value = '34.13'
$ws.Range("E50").Value = '  -5.20%  '
$ws.Range("D51").Value = '0.05720'
$ws.Range("E51").Value = '  -2.74%  '
